# Replace the paragraph about extending the system with analytics text,
# splitting the new content into two runs (as Word naturally does when the
# text is typed/edited in two bursts) plus a trailing empty run.

$d = $word.ActiveDocument

$oldText = "Ukoliko bi naš sistem bio proširen sa korisničkim zahtjevom da aplikacija prati rezultate korisnika studenata, kao dodatni kriterij kvalitete tutora i instrukcija bilo bi neophodno prosiriti i bazu podataka. To bi znatno oslabilo performance pristupa bazi podataka ukoliko bi bilo potrebno obrađivati podatke za potrebe statistike ili pretrage optimalnog tutora po željenim kriterijima. Zbog toga je korisno da se implementira Prototype pattern koji bi omogućio da već učitane podatke "

$newPart1 = "Ukoliko bi željeli da analiziramo karakteristike određenih tutora i studenata, što bi znatno oslabilo performance dobavljanja informacija iz baze podataka i samim time usporilo rad naše aplikacije. "
$newPart2 = " Zbog toga je korisno da se implementira Prototype pattern koji bi omogućio da već učitane podatke “recikliramo”, odnosno ponovno iskoristimo putem kloniranja."

# Locate the target paragraph by scanning for the unique old text.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*$oldText*") {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph"
}

$full = $target.Range
$s = $full.Start
$e = $full.End
# Exclude the trailing paragraph mark from the replaceable text range.
$textEnd = $e - 1

# 1) Replace the whole paragraph's text (minus paragraph mark) with the
#    concatenation of the two new sentences.
$r1 = $d.Range($s, $textEnd)
$r1.Text = $newPart1 + $newPart2

# 2) Split the merged run into two runs at the boundary between the two
#    sentences by toggling a formatting property off/on (this forces a run
#    break without altering visible formatting).
$splitPos = $s + $newPart1.Length
$splitRange = $d.Range($splitPos, $r1.End)
$splitRange.Bold = 1
$splitRange.Bold = 0

# 3) Create a trailing empty run (matching the target structure) by
#    temporarily inserting a paragraph break right after the text, then
#    removing that break again -- this leaves behind an empty run with the
#    same formatting instead of re-merging into run 2.
$insertPoint = $d.Range($r1.End, $r1.End)
$insertPoint.InsertParagraphAfter()

$markRange = $d.Range($r1.End, $r1.End + 1)
$markRange.Delete()

Write-Output "Replacement complete"
